# New crime data collected - weekly CompStat refresh (66th Precinct)
# Updates: report header (volume / week-of dates), column E width, and the
# weekly/28-day/YTD/2-year crime figures for rows 16-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Write a plain number into a cell, leaving its existing style/number format
# alone (use when the cell was already numeric before and after the edit).
function Set-Num($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# Convert a cell to a literal (shared-string) text value while reusing the
# exact cell style of a same-styled donor cell elsewhere on the sheet, so the
# "text" cell look (right aligned, General number format) matches cells like
# C14/C23 instead of picking up a brand new style. The leading apostrophe
# forces Excel to store the value as text instead of re-parsing it as a
# number (important for values like "0").
function Set-TextCell($ws, $addr, $donorAddr, $text) {
    $dst = $ws.Range($addr)
    $dst.Value = "'" + $text
    $ws.Range($donorAddr).Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Convert a cell that currently holds text back into a real number, reusing
# the style of a same-styled numeric donor cell.
function Set-NumCell($ws, $addr, $donorAddr, $val) {
    $dst = $ws.Range($addr)
    $dst.Value = $val
    $ws.Range($donorAddr).Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Report header: Volume/Number and the week-of date range
# ---------------------------------------------------------------------------

$ws.Range("A8").Value = "Volume 31   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"

# ---------------------------------------------------------------------------
# Column E got narrower (matches column D/F's width now)
# ---------------------------------------------------------------------------

$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# ---------------------------------------------------------------------------
# Row 16 - Murder
# ---------------------------------------------------------------------------

Set-Num $ws "F16" 12
Set-Num $ws "H16" 71.428571428571
Set-Num $ws "I16" 57
Set-Num $ws "J16" 42
Set-Num $ws "K16" 35.714285714285
Set-Num $ws "L16" 35.714285714285
Set-Num $ws "M16" -30.487804878048
Set-Num $ws "N16" -85

# ---------------------------------------------------------------------------
# Row 17 - Rape
# ---------------------------------------------------------------------------

Set-Num $ws "C17" 3
Set-Num $ws "D17" 1
Set-Num $ws "E17" 200
Set-Num $ws "F17" 14
Set-Num $ws "G17" 18
Set-Num $ws "H17" -22.222222222222
Set-Num $ws "I17" 111
Set-Num $ws "J17" 110
Set-Num $ws "K17" 0.909090909090
Set-Num $ws "L17" 2.777777777777
Set-Num $ws "M17" 46.052631578947
Set-Num $ws "N17" -31.901840490797

# ---------------------------------------------------------------------------
# Row 18 - Robbery
# ---------------------------------------------------------------------------

Set-Num $ws "C18" 2
Set-Num $ws "D18" 3
Set-Num $ws "E18" -33.333333333333
Set-Num $ws "F18" 9
Set-Num $ws "G18" 11
Set-Num $ws "H18" -18.181818181818
Set-Num $ws "I18" 58
Set-Num $ws "J18" 62
Set-Num $ws "K18" -6.451612903225
Set-Num $ws "L18" -32.558139534883
Set-Num $ws "M18" -70.408163265306
Set-Num $ws "N18" -93.103448275862

# ---------------------------------------------------------------------------
# Row 19 - Fel. Assault
# ---------------------------------------------------------------------------

Set-Num $ws "C19" 10
Set-Num $ws "D19" 3
Set-Num $ws "E19" 233.333333333333
Set-Num $ws "F19" 41
Set-Num $ws "G19" 44
Set-Num $ws "H19" -6.818181818181
Set-Num $ws "I19" 281
Set-Num $ws "J19" 283
Set-Num $ws "K19" -0.706713780918
Set-Num $ws "L19" -8.469055374592
Set-Num $ws "M19" 38.423645320197
Set-Num $ws "N19" -10.793650793650

# ---------------------------------------------------------------------------
# Row 20 - Burglary
# ---------------------------------------------------------------------------

Set-Num $ws "C20" 3
Set-Num $ws "D20" 6
Set-Num $ws "E20" -50
Set-Num $ws "F20" 14
Set-Num $ws "G20" 18
Set-Num $ws "H20" -22.222222222222
Set-Num $ws "I20" 92
Set-Num $ws "J20" 82
Set-Num $ws "K20" 12.195121951219
Set-Num $ws "L20" 91.666666666666
Set-Num $ws "M20" 26.027397260274
Set-Num $ws "N20" -90.631364562118

# ---------------------------------------------------------------------------
# Row 21 - Gr. Larceny (TOTAL-styled row)
# ---------------------------------------------------------------------------

Set-Num $ws "C21" 20
Set-Num $ws "D21" 15
Set-Num $ws "E21" 33.333333333333
Set-Num $ws "F21" 90
Set-Num $ws "G21" 99
Set-Num $ws "H21" -9.090909090909
Set-Num $ws "I21" 609
Set-Num $ws "J21" 591
Set-Num $ws "K21" 3.045685279187
Set-Num $ws "L21" 0.995024875621
Set-Num $ws "M21" -3.639240506329
Set-Num $ws "N21" -77.410979228486

# ---------------------------------------------------------------------------
# Row 22 - G.L.A. (Week-to-date columns C/D/E flip to "no data" text)
# ---------------------------------------------------------------------------

Set-TextCell $ws "C22" "C14" "0"
Set-Num $ws "F22" 1
Set-Num $ws "H22" 0

# ---------------------------------------------------------------------------
# Row 24 - Transit
# ---------------------------------------------------------------------------

Set-Num $ws "C24" 18
Set-Num $ws "D24" 16
Set-Num $ws "E24" 12.5
Set-Num $ws "F24" 80
Set-Num $ws "G24" 91
Set-Num $ws "H24" -12.087912087912
Set-Num $ws "I24" 512
Set-Num $ws "J24" 584
Set-Num $ws "K24" -12.328767123287
Set-Num $ws "L24" -11.111111111111
Set-Num $ws "M24" 12.527472527472

# ---------------------------------------------------------------------------
# Row 25 - Housing (Week-to-date column C flips from "no data" back to real data)
# ---------------------------------------------------------------------------

Set-NumCell $ws "C25" "C26" 4
Set-Num $ws "D25" 8
Set-Num $ws "E25" -50
Set-Num $ws "G25" 19
Set-Num $ws "H25" -15.789473684210
Set-Num $ws "I25" 123
Set-Num $ws "J25" 168
Set-Num $ws "K25" -26.785714285714
Set-Num $ws "L25" -28.488372093023

# ---------------------------------------------------------------------------
# Row 26 - Petit Larceny
# ---------------------------------------------------------------------------

Set-Num $ws "C26" 7
Set-Num $ws "D26" 10
Set-Num $ws "E26" -30
Set-Num $ws "F26" 32
Set-Num $ws "G26" 37
Set-Num $ws "H26" -13.513513513513
Set-Num $ws "I26" 211
Set-Num $ws "J26" 178
Set-Num $ws "K26" 18.539325842696
Set-Num $ws "L26" 27.108433734939
Set-Num $ws "M26" -10.593220338983

# ---------------------------------------------------------------------------
# Row 28 - Misd. Assault
# ---------------------------------------------------------------------------

Set-Num $ws "C28" 1
Set-Num $ws "D28" 1
Set-Num $ws "E28" 0
Set-Num $ws "F28" 7
Set-Num $ws "G28" 7
Set-Num $ws "H28" 0
Set-Num $ws "I28" 38
Set-Num $ws "J28" 33
Set-Num $ws "K28" 15.151515151515
Set-Num $ws "L28" -5

# ---------------------------------------------------------------------------
# Row 31 - Shooting Inc. (Week-to-date columns D/E flip to "no data" text)
# ---------------------------------------------------------------------------

Set-TextCell $ws "D31" "D23" "0"
Set-TextCell $ws "E31" "E23" "***.*"

# ---------------------------------------------------------------------------
# Row 33 - Traffic Fatalities (Week-to-date columns D/E/F reshuffle)
# ---------------------------------------------------------------------------

Set-TextCell $ws "D33" "D23" "0"
Set-TextCell $ws "E33" "E23" "***.*"
Set-TextCell $ws "F33" "F23" "0"
Set-Num $ws "G33" 1

Write-Output "CompStat weekly figures updated"
